$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last row (row 43), pushing the
# existing row 43 record down to row 45, then fill in the new row 43/44
# records and the (shifted) row 45 record.
$ws.Rows.Item(43).Resize(2).Insert()

# Row 43: Early Diamond / Segunda
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(43, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(43, 4).Value = 44568
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
$ws.Cells.Item(43, 5).Value = 15
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100103
$ws.Cells.Item(43, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(43, 9).Value = 100103006
$ws.Cells.Item(43, 10).Value = "Nectarín"
$ws.Cells.Item(43, 11).Value = "Early Diamond"
$ws.Cells.Item(43, 12).Value = "Segunda"
$ws.Cells.Item(43, 13).Value = 300
$ws.Cells.Item(43, 14).Value = 19000
$ws.Cells.Item(43, 15).Value = 20000
$ws.Cells.Item(43, 16).Value = 19500
$ws.Cells.Item(43, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(43, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(43, 19).Value = 1083
$ws.Cells.Item(43, 20).Value = 18

# Row 44: Super Queen / Primera
$ws.Cells.Item(44, 1).Value = 1
$ws.Cells.Item(44, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(44, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(44, 4).Value = 44568
$ws.Cells.Item(44, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
$ws.Cells.Item(44, 5).Value = 15
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100103
$ws.Cells.Item(44, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(44, 9).Value = 100103006
$ws.Cells.Item(44, 10).Value = "Nectarín"
$ws.Cells.Item(44, 11).Value = "Super Queen"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 250
$ws.Cells.Item(44, 14).Value = 20000
$ws.Cells.Item(44, 15).Value = 22000
$ws.Cells.Item(44, 16).Value = 21000
$ws.Cells.Item(44, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(44, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(44, 19).Value = 1167
$ws.Cells.Item(44, 20).Value = 18

# Row 45: previously row 43 (Artic Pride / Segunda) — values unchanged, just shifted down
$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(45, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(45, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(45, 4).Value = 44540
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
$ws.Cells.Item(45, 5).Value = 15
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100103
$ws.Cells.Item(45, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(45, 9).Value = 100103006
$ws.Cells.Item(45, 10).Value = "Nectarín"
$ws.Cells.Item(45, 11).Value = "Artic Pride"
$ws.Cells.Item(45, 12).Value = "Segunda"
$ws.Cells.Item(45, 13).Value = 250
$ws.Cells.Item(45, 14).Value = 21000
$ws.Cells.Item(45, 15).Value = 22000
$ws.Cells.Item(45, 16).Value = 21500
$ws.Cells.Item(45, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(45, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(45, 19).Value = 1194
$ws.Cells.Item(45, 20).Value = 18
